$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $value
    $cell.ClearFormats()
}

Set-TextValue "D2" "63.342.37"
Set-TextValue "E2" "  -7.24%  "
Set-TextValue "D3" "3.273.60"
Set-TextValue "E3" "  -8.44%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "181.22"
Set-TextValue "E5" "  -12.06%  "
Set-TextValue "D6" "515.99"
Set-TextValue "E6" "  -7.81%  "
Set-TextValue "E7" "  -2.00%  "
Set-TextValue "D8" "3.270.93"
Set-TextValue "E8" "  -8.41%  "
Set-TextValue "E9" "  +0.04%  "
Set-TextValue "E10" "  -8.39%  "
Set-TextValue "D11" "58.76"
Set-TextValue "E11" "  -6.96%  "
Set-TextValue "E12" "  -10.26%  "
Set-TextValue "D13" "0.0000254"
Set-TextValue "E13" "  -8.72%  "
Set-TextValue "D14" "9.10"
Set-TextValue "E14" "  -9.65%  "
Set-TextValue "D15" "3.803.38"
Set-TextValue "E15" "  -8.46%  "
Set-TextValue "E16" "  -5.81%  "
Set-TextValue "D17" "3.283.83"
Set-TextValue "E17" "  -8.14%  "
Set-TextValue "D18" "17.56"
Set-TextValue "E18" "  -7.92%  "
Set-TextValue "D19" "63.330.07"
Set-TextValue "E19" "  -7.09%  "
Set-TextValue "D20" "10.93"
Set-TextValue "E20" "  -9.73%  "
Set-TextValue "D21" "0.944"
Set-TextValue "E21" "  -10.69%  "
Set-TextValue "D22" "370.23"
Set-TextValue "E22" "  -7.67%  "
Set-TextValue "D23" "11.17"
Set-TextValue "E23" "  -8.41%  "
Set-TextValue "D24" "80.01"
Set-TextValue "E24" "  -5.13%  "
Set-TextValue "D25" "3.65"
Set-TextValue "E25" "  -11.11%  "
Set-TextValue "D26" "3.85"
Set-TextValue "E26" "  +0.80%  "
Set-TextValue "D27" "5.97"
Set-TextValue "E27" "  -2.59%  "
Set-TextValue "D28" "2.64"
Set-TextValue "E28" "  -7.82%  "
Set-TextValue "D29" "11.33"
Set-TextValue "E29" "  -8.44%  "
Set-TextValue "D30" "8.28"
Set-TextValue "E30" "  -8.35%  "
Set-TextValue "B31" "EthereumClassic"
Set-TextValue "C31" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D31" "28.40"
Set-TextValue "E31" "  -9.17%  "
Set-TextValue "B32" "Bittensor"
Set-TextValue "C32" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D32" "644.52"
Set-TextValue "E32" "  -10.58%  "
Set-TextValue "D33" "6.66"
Set-TextValue "E33" "  -11.11%  "
Set-TextValue "D34" "11.14"
Set-TextValue "E34" "  -7.15%  "
Set-TextValue "D35" "59.32"
Set-TextValue "E35" "  -7.10%  "
Set-TextValue "D36" "0.104"
Set-TextValue "E36" "  -6.46%  "
Set-TextValue "E37" "  -0.19%  "
Set-TextValue "D38" "0.385"
Set-TextValue "E38" "  -8.01%  "
Set-TextValue "D39" "35.90"
Set-TextValue "E39" "  -12.43%  "
Set-TextValue "D40" "0.997"
Set-TextValue "E40" "  -0.25%  "
Set-TextValue "D41" "2.965.99"
Set-TextValue "E41" "  -6.24%  "
Set-TextValue "D42" "0.124"
Set-TextValue "E42" "  -5.79%  "
Set-TextValue "D43" "0.0₃0647"
Set-TextValue "E43" "  -10.53%  "
Set-TextValue "D44" "2.67"
Set-TextValue "E44" "  -17.25%  "
Set-TextValue "D45" "2.41"
Set-TextValue "E45" "  -6.12%  "
Set-TextValue "D46" "2.58"
Set-TextValue "E46" "  -6.13%  "
Set-TextValue "D47" "0.0386"
Set-TextValue "E47" "  -5.87%  "
Set-TextValue "D48" "2.78"
Set-TextValue "E48" "  +3.58%  "
Set-TextValue "B49" "ApeXProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D49" "2.97"
Set-TextValue "E49" "  -4.00%  "
Set-TextValue "B50" "Stellar"
Set-TextValue "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.125"
Set-TextValue "E50" "  -3.95%  "
Set-TextValue "B51" "dogwifhat"
Set-TextValue "C51" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "2.43"
Set-TextValue "E51" "  -21.65%  "
